$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1526.6666
$ws.Range("I43").Value = 780
$ws.Range("J43").Value = 1900
$ws.Range("K43").Value = 780
$ws.Range("L43").Value = 1900
$ws.Range("M43").Value = -711
$ws.Range("N43").Value = -2038

$ws.Range("H116").Value = 7764.316
$ws.Range("I116").Value = 8398.134
$ws.Range("J116").Value = 5387.5
$ws.Range("K116").Value = 8398.134
$ws.Range("L116").Value = 5387.5
$ws.Range("M116").Value = -4956.134
$ws.Range("N116").Value = -12271.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 21659.666
$ws.Range("I37").Value = 10000
$ws.Range("J37").Value = 27489.5
$ws.Range("K37").Value = 10000
$ws.Range("L37").Value = 27489.5
$ws.Range("M37").Value = -9727
$ws.Range("N37").Value = -28035.5

$ws.Range("H44").Value = 20000
$ws.Range("J44").Value = 20000
$ws.Range("L44").Value = 20000
$ws.Range("N44").Value = -20976

$ws.Range("H45").Value = 9391.75
$ws.Range("I45").Value = 11909.667
$ws.Range("J45").Value = 1838
$ws.Range("K45").Value = 11909.667
$ws.Range("L45").Value = 1838
$ws.Range("M45").Value = -11532.667
$ws.Range("N45").Value = -2592

$ws.Range("H55").Value = 19899.666
$ws.Range("J55").Value = 19899.666
$ws.Range("L55").Value = 19899.666
$ws.Range("N55").Value = -20529.666

$ws.Range("H61").Value = 347868.7
$ws.Range("I61").Value = 2437.2666
$ws.Range("J61").Value = 717973.8
$ws.Range("K61").Value = 2437.2666
$ws.Range("L61").Value = 717973.8
$ws.Range("M61").Value = -2225.2666
$ws.Range("N61").Value = -718397.8

$ws.Range("H63").Value = 200004600
$ws.Range("I63").Value = 200004600
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 200004600
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -200003914
$ws.Range("N63").Value = $null

$ws.Range("H66").Value = 200004600
$ws.Range("I66").Value = 200004600
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 1000023000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -1000019568
$ws.Range("N66").Value = $null

$ws.Range("H136").Value = 347868.7
$ws.Range("I136").Value = 2437.2666
$ws.Range("J136").Value = 717973.8
$ws.Range("K136").Value = 7311.7998
$ws.Range("L136").Value = 2153921.4
$ws.Range("M136").Value = -4761.7998
$ws.Range("N136").Value = -2159021.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13911.45
$ws.Range("I20").Value = 1326.8334
$ws.Range("J20").Value = 32788.375
$ws.Range("K20").Value = 1326.8334
$ws.Range("L20").Value = 32788.375
$ws.Range("M20").Value = -1079.8334
$ws.Range("N20").Value = -33282.375

$ws.Range("H105").Value = 12585.25
$ws.Range("I105").Value = 17819.77
$ws.Range("J105").Value = 2864
$ws.Range("K105").Value = 17819.77
$ws.Range("L105").Value = 2864
$ws.Range("M105").Value = -16072.77
$ws.Range("N105").Value = -6358

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1860.091
$ws.Range("I16").Value = 1519.7333
$ws.Range("K16").Value = 1519.7333
$ws.Range("M16").Value = -1232.7333

$ws.Range("H22").Value = 887.8
$ws.Range("I22").Value = 859.75
$ws.Range("K22").Value = 859.75
$ws.Range("M22").Value = -509.75

$ws.Range("H31").Value = 3102.647
$ws.Range("I31").Value = 2042.3334
$ws.Range("J31").Value = 4102.3716
$ws.Range("K31").Value = 2042.3334
$ws.Range("L31").Value = 4102.3716
$ws.Range("M31").Value = -1747.3334
$ws.Range("N31").Value = -4692.3716

$ws.Range("H34").Value = 3102.647
$ws.Range("I34").Value = 2042.3334
$ws.Range("J34").Value = 4102.3716
$ws.Range("K34").Value = 2042.3334
$ws.Range("L34").Value = 4102.3716
$ws.Range("M34").Value = -1840.3334
$ws.Range("N34").Value = -4506.3716

$ws.Range("H58").Value = 1470.9556
$ws.Range("I58").Value = 1203.2963
$ws.Range("J58").Value = 1872.4445
$ws.Range("K58").Value = 1203.2963
$ws.Range("L58").Value = 1872.4445
$ws.Range("M58").Value = -1000.2963
$ws.Range("N58").Value = -2278.4445

$ws.Range("H113").Value = 1860.091
$ws.Range("I113").Value = 1519.7333
$ws.Range("K113").Value = 1519.7333
$ws.Range("M113").Value = 650.2666999999999

$ws.Range("H134").Value = 2746.775
$ws.Range("I134").Value = 2827.7646
$ws.Range("J134").Value = 2287.8333
$ws.Range("K134").Value = 8483.293799999999
$ws.Range("L134").Value = 6863.499899999999
$ws.Range("M134").Value = -5948.293799999999
$ws.Range("N134").Value = -11933.4999

$ws.Range("H136").Value = 1470.9556
$ws.Range("I136").Value = 1203.2963
$ws.Range("J136").Value = 1872.4445
$ws.Range("K136").Value = 3609.8889
$ws.Range("L136").Value = 5617.333500000001
$ws.Range("M136").Value = -1059.8889
$ws.Range("N136").Value = -10717.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1125.25
$ws.Range("I34").Value = 900.6667
$ws.Range("J34").Value = 1260
$ws.Range("K34").Value = 2702.0001
$ws.Range("L34").Value = 3780
$ws.Range("M34").Value = -2618.0001
$ws.Range("N34").Value = -3948

$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").Value = $null

$ws.Range("H124").Value = 10000
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 10000
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 30000
$ws.Range("M124").Value = $null
$ws.Range("N124").Value = -39820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 386693.38
$ws.Range("I102").Value = 808072.7
$ws.Range("J102").Value = 1955.7391
$ws.Range("K102").Value = 808072.7
$ws.Range("L102").Value = 1955.7391
$ws.Range("M102").Value = -806450.7
$ws.Range("N102").Value = -5199.7391

$ws.Range("H113").Value = 1794.85
$ws.Range("I113").Value = 1199.7
$ws.Range("J113").Value = 2390
$ws.Range("K113").Value = 1199.7
$ws.Range("L113").Value = 2390
$ws.Range("M113").Value = 970.3
$ws.Range("N113").Value = -6730

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2483.4443
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2483.4443
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2483.4443
$ws.Range("M22").Value = $null
$ws.Range("N22").Value = -3073.4443

$ws.Range("H27").Value = 2483.4443
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2483.4443
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2483.4443
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = -2697.4443

$ws.Range("H40").Value = 2344.389
$ws.Range("I40").Value = 2449.9375
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 2449.9375
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -2313.9375
$ws.Range("N40").Value = -1772

$ws.Range("H46").Value = 977
$ws.Range("I46").Value = 783.3333
$ws.Range("J46").Value = 1122.25
$ws.Range("K46").Value = 783.3333
$ws.Range("L46").Value = 1122.25
$ws.Range("M46").Value = -595.3333
$ws.Range("N46").Value = -1498.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1202.4783
$ws.Range("I136").Value = 657
$ws.Range("J136").Value = 2225.25
$ws.Range("K136").Value = 1971
$ws.Range("L136").Value = 6675.75
$ws.Range("M136").Value = 579
$ws.Range("N136").Value = -11775.75
